$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I6").Value = 1.4
$ws.Range("J6").Value = 7
$ws.Range("K6").Value = 2.3
$ws.Range("L6").Value = 1.91
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.67
$ws.Range("W6").Value = 17
$ws.Range("Y6").Value = 23
$ws.Range("AC6").Value = 11
$ws.Range("AE6").Value = 21
$ws.Range("AF6").Value = 67
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 7
$ws.Range("AI6").Value = 6.5
$ws.Range("AL6").Value = 13
$ws.Range("AM6").Value = 29
$ws.Range("AO6").Value = 41
$ws.Range("AQ6").Value = 151
$ws.Range("AU6").Value = 9.5
$ws.Range("AV6").Value = 67
$ws.Range("AY6").Value = 19
$ws.Range("BB6").Value = 151
